$wb = $excel.ActiveWorkbook

# Sheet "展览" - update 想去人数 (attendee counts) for several rows
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 1123
$wsExhibit.Range("F4").Value = 1823
$wsExhibit.Range("F5").Value = 803
$wsExhibit.Range("F6").Value = 464

# Sheet "全部类型" - same underlying rows, updated counts
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 1123
$wsAll.Range("F4").Value = 1823
$wsAll.Range("F6").Value = 803
$wsAll.Range("F7").Value = 464
